$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.871.41"
$ws.Range("E2").Value = "  -0.85%  "

$ws.Range("D3").Value = "2.500.25"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").Value = "'537.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("D6").Value = "'136.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.61%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.71%  "

$ws.Range("D9").Value = "2.526.61"
$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("D11").Value = "'0.157"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.19%  "

$ws.Range("D12").Value = "'5.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.23%  "

$ws.Range("E13").Value = "  -1.52%  "

$ws.Range("D14").Value = "2.956.34"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").Value = "'23.04"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "58.817.24"
$ws.Range("E16").Value = "  -0.98%  "

$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").Value = "2.525.97"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("D19").Value = "'11.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").Value = "'323.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.54%  "

$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").Value = "'5.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.44%  "

$ws.Range("D24").Value = "'65.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.36%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("E26").Value = "  -1.93%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").Value = "'7.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.59%  "

$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'6.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.72%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0774"
$ws.Range("E30").Value = "  -0.56%  "

$ws.Range("E31").Value = "  -1.34%  "

$ws.Range("D32").Value = "'167.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.44%  "

$ws.Range("E33").Value = "  +5.56%  "

$ws.Range("E35").Value = "  +1.58%  "

$ws.Range("D36").Value = "'18.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.34%  "

$ws.Range("E37").Value = "  -2.61%  "

$ws.Range("E38").Value = "  -3.28%  "

$ws.Range("D39").Value = "'36.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.70%  "

$ws.Range("D40").Value = "'0.814"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.37%  "

$ws.Range("E41").Value = "  -1.57%  "

$ws.Range("D42").Value = "'284.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.17%  "

$ws.Range("D43").Value = "'5.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.98%  "

$ws.Range("D44").Value = "'132.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.43%  "

$ws.Range("D45").Value = "'0.994"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.39%  "

$ws.Range("D46").Value = "'0.606"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.46%  "

$ws.Range("D47").Value = "'10.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.28%  "

$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("D50").Value = "'0.0220"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.70%  "

$ws.Range("D51").Value = "'17.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.67%  "
